$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.123.13'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.401.53'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.50'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.75'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range('E8').Value = '  +1.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.409.25'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.108'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.18%  '
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.21'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.41%  '
$ws.Range('E13').Value = '  +2.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.46'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.10%  '
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.838.48'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.996.26'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.410.21'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.07'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.14%  '
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.45'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.07'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.52%  '
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.93'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +5.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '65.12'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '592.30'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.23'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0946'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.38%  '
$ws.Range('E30').Value = '  -0.96%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.00'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.37'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.67%  '
$ws.Range('E33').Value = '  -0.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.132'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.69%  '
$ws.Range('E35').Value = '  +3.42%  '
$ws.Range('E36').Value = '  -0.26%  '
$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.372'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.11%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.62'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '152.20'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.26'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.09%  '
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('E43').Value = '  +1.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.77'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.49'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +5.64%  '
$ws.Range('E46').Value = '  +2.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '141.60'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.16%  '
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.591'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.68%  '
$ws.Range('E50').Value = '  +1.16%  '
